# Applies the "January 2020" monthly update:
#  - collapses the stale B4:C50 review-selections on the earlier month
#    sheets down to the single cell the user was last on
#  - duplicates "December 2019" into a new "January 2020" sheet
#  - refreshes the duplicated sheet's header label + rank/name/contribution
#    rows with next month's data
#  - leaves "January 2020" as the active/selected sheet

$wb = $excel.ActiveWorkbook

# ---- collapse leftover multi-cell selections on the existing sheets ----
$aug = $wb.Worksheets.Item("August 2019")
$aug.Activate() | Out-Null
$aug.Range("C49").Select() | Out-Null

$sep = $wb.Worksheets.Item("September 2019")
$sep.Activate() | Out-Null
$sep.Range("C15").Select() | Out-Null

$oct = $wb.Worksheets.Item("October 2019")
$oct.Activate() | Out-Null
$oct.Range("E3").Select() | Out-Null

$nov = $wb.Worksheets.Item("November 2019")
$nov.Activate() | Out-Null
$nov.Range("E2").Select() | Out-Null

$dec = $wb.Worksheets.Item("December 2019")
$dec.Activate() | Out-Null
$dec.Range("E24").Select() | Out-Null

# ---- duplicate December 2019 -> January 2020 (placed right after it) ----
$dec.Copy($null, $dec) | Out-Null
$jan = $wb.Worksheets.Item($wb.Worksheets.Count)
$jan.Name = "January 2020"

# fix the JSON-ish header label formula that was copied verbatim
$jan.Range("E2").Formula = '="    """&"January 2020"&""""&":"'

# ---- next month's leaderboard ----
$ranks = 1..50
$names = @(
  "Eternal","Smile","Savages","Bounce","Elite","Spring","Sunset","Epic","Beaters","Downtime",
  "lolicafe","RainSong","Imperium","Remorse","Undertale","Gintama","Maha","Erda","Cleanse","Tama",
  "Broke","Atelier","Sora","Lithe","Revive","Oceania","Rising","Sugar","Ravers","Artifacts",
  "Fabled","Earnest","Aloe","Skyfall","Mystical","CyberThreat","Howl","Fandom","chigga","RainDrop",
  "Path","Kingdom","Exorcist","Bubbles","Coffee","Weibo","Faction","Reboot","HeavenSent","Comity"
)
$contributions = @(
  243043578,237838196,228672705,212909088,207527311,147106317,138924603,133258625,125747064,120319160,
  102061033,99680946,99350578,89190116,88305254,87768163,84739973,84517975,82608759,79091632,
  78951041,78316757,75722680,73552585,73029232,68152097,65628617,63988479,63310658,61627857,
  56716634,56608216,56444050,52169401,50864570,50519707,50488967,47991738,47869266,46745943,
  45219390,45107275,42688203,42478754,41798171,40867146,40254536,39616142,39129544,39056996
)

for ($i = 0; $i -lt 50; $i++) {
    $row = 4 + $i
    $jan.Cells.Item($row, 2).Value = $ranks[$i]
    $jan.Cells.Item($row, 3).Value = $names[$i]
    $jan.Cells.Item($row, 4).Value = $contributions[$i]
}

# ---- leave January 2020 selected on the cell the data entry ended at ----
$jan.Activate() | Out-Null
$jan.Range("C4").Select() | Out-Null

Write-Output "Added January 2020 sheet after December 2019."
